$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task_settings")

$ws.Range("A1").Value = "TRIMAZKON_task_1"
$ws.Range("B1").Value = "C:/Users/jakub.hlavacek.local/Desktop/JHV/test_images/Keyence/_503_Witte/datumovka/A/Height_test/"

$ws.Range("C1").Value = "'336"
$ws.Range("C1").ClearFormats()

$ws.Range("D1").Value = "'998"
$ws.Range("D1").ClearFormats()

$ws.Range("E1").Value = "12:00"
$ws.Range("F1").Value = "30.01.2025 10:03:29"

$ws.Range("G1").Value = "'"
$ws.Range("G1").ClearFormats()

